# [CC] correction of the reference values
#
# Rewrites the fmod/R reference table (A3:B23 in the original file) with
# corrected measurement values and extends it with 7 additional rows
# (A24:B30) that were missing from the original data set.
#
# Values are entered with a leading apostrophe so Excel stores them as text
# (matching the original file, where every numeric-looking value lives in
# xl/sharedStrings.xml as a <t> string, not a numeric <v>), then the
# quote-prefix style that introduces is cleared back to Normal so the cells
# keep the original's plain (style-less) look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (fmod) values, rows 3-30 (written first so the rebuilt shared-string table lists them before column B) ---
$ws.Range("A3").Value = "'15.561241"
$ws.Range("A4").Value = "'22.401348"
$ws.Range("A5").Value = "'29.266758"
$ws.Range("A6").Value = "'36.86032"
$ws.Range("A7").Value = "'38.84244"
$ws.Range("A8").Value = "'43.589268"
$ws.Range("A9").Value = "'47.646236"
$ws.Range("A10").Value = "'51.401146"
$ws.Range("A11").Value = "'56.619133"
$ws.Range("A12").Value = "'59.81015"
$ws.Range("A13").Value = "'63.843525"
$ws.Range("A14").Value = "'67.08531"
$ws.Range("A15").Value = "'70.85759"
$ws.Range("A16").Value = "'74.06288"
$ws.Range("A17").Value = "'77.41061"
$ws.Range("A18").Value = "'81.33141"
$ws.Range("A19").Value = "'87.916466"
$ws.Range("A20").Value = "'97.7758"
$ws.Range("A21").Value = "'100.33316"
$ws.Range("A22").Value = "'115.681206"
$ws.Range("A23").Value = "'119.01259"
$ws.Range("A24").Value = "'134.40422"
$ws.Range("A25").Value = "'148.2887"
$ws.Range("A26").Value = "'177.72102"
$ws.Range("A27").Value = "'205.42125"
$ws.Range("A28").Value = "'235.5931"
$ws.Range("A29").Value = "'272.31946"
$ws.Range("A30").Value = "'296.5826"

# --- Column B (R) values, rows 3-30 ---
$ws.Range("B3").Value = "'0.109677784"
$ws.Range("B4").Value = "'0.18588823"
$ws.Range("B5").Value = "'0.27217966"
$ws.Range("B6").Value = "'0.3822299"
$ws.Range("B7").Value = "'0.40800852"
$ws.Range("B8").Value = "'0.4784683"
$ws.Range("B9").Value = "'0.52834743"
$ws.Range("B10").Value = "'0.5683734"
$ws.Range("B11").Value = "'0.6097778"
$ws.Range("B12").Value = "'0.6274577"
$ws.Range("B13").Value = "'0.6473149"
$ws.Range("B14").Value = "'0.6522843"
$ws.Range("B15").Value = "'0.6521323"
$ws.Range("B16").Value = "'0.6452131"
$ws.Range("B17").Value = "'0.6333704"
$ws.Range("B18").Value = "'0.62010473"
$ws.Range("B19").Value = "'0.57756937"
$ws.Range("B20").Value = "'0.51043147"
$ws.Range("B21").Value = "'0.48305747"
$ws.Range("B22").Value = "'0.38639748"
$ws.Range("B23").Value = "'0.36280885"
$ws.Range("B24").Value = "'0.29872113"
$ws.Range("B25").Value = "'0.2538336"
$ws.Range("B26").Value = "'0.18619592"
$ws.Range("B27").Value = "'0.1454662"
$ws.Range("B28").Value = "'0.11394783"
$ws.Range("B29").Value = "'0.089489736"
$ws.Range("B30").Value = "'0.07725139"

# Strip the quote-prefix style that typing a leading apostrophe introduces, so the
# corrected cells keep the plain (no explicit style) look of the original file.
$ws.Range("A3:B30").Style = "Normal"

# Selection / scroll state to mirror the author reopening the sheet scrolled to the table.
$ws.Range("B3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
